# Update "想去人数" (number of people interested) counts in column F
# for the sheets "展览" and "全部类型", reflecting newly refreshed data.

$wb = $excel.ActiveWorkbook

# Row -> new value mapping shared by both sheets.
$updates = @{
    2  = 194
    4  = 159
    5  = 1321
    6  = 18362
    8  = 267
    10 = 6900
    11 = 693
    12 = 163
    14 = 118
    19 = 262
    22 = 38
    25 = 280
    28 = 5176
    29 = 541
    33 = 12139
    34 = 1292
    38 = 3932
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}

# The "全部类型" sheet contains the same rows, but shifted down by two
# (rows 30/31 hold extra entries from the "演出" sheet), so rows >= 33
# in "展览" correspond to rows >= 35 in "全部类型".
$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    if ($row -ge 33) {
        $targetRow = $row + 2
    } else {
        $targetRow = $row
    }
    $ws4.Cells.Item($targetRow, 6).Value = $updates[$row]
}
